$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 44595
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112038
$ws.Cells.Item($row, 7).Value = "Cebollín baby"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 250
$ws.Cells.Item($row, 11).Value = 3500
$ws.Cells.Item($row, 12).Value = 4000
$ws.Cells.Item($row, 13).Value = 3750
$ws.Cells.Item($row, 14).Value = '$/paquete 1,5 a 2 kilos'
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 1875
$ws.Cells.Item($row, 17).Value = 2
$ws.Cells.Item($row, 18).Value = "Hortaliza"
